$d = $word.ActiveDocument

# The document has a Pearson logo image (image1.png) embedded in both
# footers, and a BTEC logo image (image2.jpg) embedded in both headers.
# The commit swaps each image's shape-name metadata with the other
# image's old file-name-derived name:
#   footers: "image1.png" -> "image2.png"
#   headers: "image2.jpg" -> "image1.jpg"
#
# Headers: a direct property-chain assignment works reliably.
$d.Sections(1).Headers(1).Range.InlineShapes(1).Name = "image1.jpg"
$d.Sections(1).Headers(2).Range.InlineShapes(1).Name = "image1.jpg"

# Footers: the direct chain leaves a stale handle on this engine, so
# select the inline picture first and rename it through $word.Selection.
$ftr1 = $d.Sections(1).Footers(1)
$ftr1Shape = $ftr1.Range.InlineShapes(1)
[void]$ftr1Shape.Select()
$word.Selection.InlineShapes(1).Name = "image2.png"

$ftr2 = $d.Sections(1).Footers(2)
$ftr2Shape = $ftr2.Range.InlineShapes(1)
[void]$ftr2Shape.Select()
$word.Selection.InlineShapes(1).Name = "image2.png"
